$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 24923
$ws.Range("E2").Value = 2734
$ws.Range("F2").Value = 2734
$ws.Range("G2").Value = 2500
$ws.Range("H2").Value = 1861
$ws.Range("I2").Value = 1566
$ws.Range("J2").Value = 294
$ws.Range("K2").Value = 79314
$ws.Range("L2").Value = 43624
$ws.Range("M2").Value = 35691
$ws.Range("N2").Value = 26157
$ws.Range("O2").Value = 9533
$ws.Range("P2").Value = 492
$ws.Range("Q2").Value = 3787
$ws.Range("R2").Value = -5139
$ws.Range("S2").Value = 1855
$ws.Range("T2").Value = 4164
$ws.Range("U2").Value = -377
$ws.Range("V2").Value = 24369
$ws.Range("W2").Value = 10.97
$ws.Range("X2").Value = 7.46
$ws.Range("Y2").Value = 6.23
$ws.Range("Z2").Value = 2.41
$ws.Range("AA2").Value = 122.23
$ws.Range("AB2").Value = 3896.9
$ws.Range("AC2").Value = 15908
$ws.Range("AD2").Value = 11.38
$ws.Range("AE2").Value = 265924
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 1150
$ws.Range("AH2").Value = 0.64
$ws.Range("AI2").Value = 7.22
$ws.Range("AJ2").Value = 9845181
$ws.Range("D3").Value = 25640
$ws.Range("E3").Value = 2621
$ws.Range("F3").Value = 2621
$ws.Range("G3").Value = 5834
$ws.Range("H3").Value = 4332
$ws.Range("I3").Value = 4021
$ws.Range("J3").Value = 311
$ws.Range("K3").Value = 79182
$ws.Range("L3").Value = 38573
$ws.Range("M3").Value = 40609
$ws.Range("N3").Value = 30233
$ws.Range("O3").Value = 10376
$ws.Range("P3").Value = 492
$ws.Range("Q3").Value = 2786
$ws.Range("R3").Value = -1508
$ws.Range("S3").Value = -1520
$ws.Range("T3").Value = 6199
$ws.Range("U3").Value = -3413
$ws.Range("V3").Value = 19184
$ws.Range("W3").Value = 10.22
$ws.Range("X3").Value = 16.89
$ws.Range("Y3").Value = 14.26
$ws.Range("Z3").Value = 5.47
$ws.Range("AA3").Value = 94.98999999999999
$ws.Range("AB3").Value = 4670.01
$ws.Range("AC3").Value = 40843
$ws.Range("AD3").Value = 5.63
$ws.Range("AE3").Value = 307358
$ws.Range("AF3").Value = 0.75
$ws.Range("AG3").Value = 1250
$ws.Range("AH3").Value = 0.54
$ws.Range("AI3").Value = 3.06
$ws.Range("AJ3").Value = 9845181
$ws.Range("D4").Value = 29475
$ws.Range("E4").Value = 2514
$ws.Range("F4").Value = 2514
$ws.Range("G4").Value = 4154
$ws.Range("H4").Value = 3234
$ws.Range("I4").Value = 2271
$ws.Range("J4").Value = 962
$ws.Range("K4").Value = 99545
$ws.Range("L4").Value = 52616
$ws.Range("M4").Value = 46929
$ws.Range("N4").Value = 32341
$ws.Range("O4").Value = 14588
$ws.Range("P4").Value = 492
$ws.Range("Q4").Value = 3594
$ws.Range("R4").Value = -8956
$ws.Range("S4").Value = 5703
$ws.Range("T4").Value = 7926
$ws.Range("U4").Value = -4332
$ws.Range("V4").Value = 26529
$ws.Range("W4").Value = 8.529999999999999
$ws.Range("X4").Value = 10.97
$ws.Range("Y4").Value = 7.26
$ws.Range("Z4").Value = 3.62
$ws.Range("AA4").Value = 112.12
$ws.Range("AB4").Value = 5078.38
$ws.Range("AC4").Value = 23071
$ws.Range("AD4").Value = 7.63
$ws.Range("AE4").Value = 328792
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 0.71
$ws.Range("AI4").Value = 5.41
$ws.Range("AJ4").Value = 9845181
$ws.Range("D5").Value = 38714
$ws.Range("E5").Value = 3457
$ws.Range("F5").Value = 3457
$ws.Range("G5").Value = 2914
$ws.Range("H5").Value = 2136
$ws.Range("I5").Value = 1823
$ws.Range("J5").Value = 313
$ws.Range("K5").Value = 101517
$ws.Range("L5").Value = 52322
$ws.Range("M5").Value = 49195
$ws.Range("N5").Value = 34218
$ws.Range("O5").Value = 14978
$ws.Range("P5").Value = 492
$ws.Range("Q5").Value = 3899
$ws.Range("R5").Value = -3048
$ws.Range("S5").Value = 120
$ws.Range("T5").Value = 3956
$ws.Range("U5").Value = -57
$ws.Range("V5").Value = 26616
$ws.Range("W5").Value = 8.93
$ws.Range("X5").Value = 5.52
$ws.Range("Y5").Value = 5.48
$ws.Range("Z5").Value = 2.13
$ws.Range("AA5").Value = 106.36
$ws.Range("AB5").Value = 5374.13
$ws.Range("AC5").Value = 18513
$ws.Range("AD5").Value = 16.2
$ws.Range("AE5").Value = 347867
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 0.42
$ws.Range("AI5").Value = 6.75
$ws.Range("AJ5").Value = 9845181
$ws.Range("D6").Value = 51857
$ws.Range("E6").Value = 3974
$ws.Range("F6").Value = 3974
$ws.Range("G6").Value = 3548
$ws.Range("H6").Value = 2849
$ws.Range("I6").Value = 2390
$ws.Range("K6").Value = 111013
$ws.Range("L6").Value = 60991
$ws.Range("M6").Value = 50022
$ws.Range("N6").Value = 34691
$ws.Range("P6").Value = 492
$ws.Range("Q6").Value = 3313
$ws.Range("R6").Value = -5926
$ws.Range("S6").Value = 4383
$ws.Range("T6").Value = 3887
$ws.Range("U6").Value = -574
$ws.Range("V6").Value = 32803
$ws.Range("W6").Value = 7.66
$ws.Range("X6").Value = 5.49
$ws.Range("Y6").Value = 6.94
$ws.Range("Z6").Value = 2.68
$ws.Range("AA6").Value = 121.93
$ws.Range("AB6").Value = 5762.54
$ws.Range("AC6").Value = 24274
$ws.Range("AD6").Value = 10.55
$ws.Range("AE6").Value = 352681
$ws.Range("AF6").Value = 0.73
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 0.78
$ws.Range("AI6").Value = 8.23
$ws.Range("AJ6").Value = 9845181
$ws.Range("D7").Value = 63678
$ws.Range("E7").Value = 4283
$ws.Range("G7").Value = 11464
$ws.Range("H7").Value = 8664
$ws.Range("I7").Value = 7969
$ws.Range("K7").Value = 135674
$ws.Range("L7").Value = 77943
$ws.Range("M7").Value = 57730
$ws.Range("N7").Value = 42014
$ws.Range("P7").Value = 491
$ws.Range("Q7").Value = 8884
$ws.Range("R7").Value = -8993
$ws.Range("S7").Value = -215
$ws.Range("T7").Value = 2429
$ws.Range("U7").Value = 7030
$ws.Range("W7").Value = 6.73
$ws.Range("X7").Value = 13.61
$ws.Range("Y7").Value = 20.78
$ws.Range("Z7").Value = 7.02
$ws.Range("AA7").Value = 135.01
$ws.Range("AC7").Value = 80943
$ws.Range("AD7").Value = 3.3
$ws.Range("AE7").Value = 427123
$ws.Range("AF7").Value = 0.63
$ws.Range("AG7").Value = 2011
$ws.Range("AH7").Value = 0.75
$ws.Range("AI7").Value = 2.48
$ws.Range("D8").Value = 69566
$ws.Range("E8").Value = 4995
$ws.Range("G8").Value = 4256
$ws.Range("H8").Value = 3212
$ws.Range("I8").Value = 2697
$ws.Range("K8").Value = 138027
$ws.Range("L8").Value = 77392
$ws.Range("M8").Value = 60635
$ws.Range("N8").Value = 44496
$ws.Range("P8").Value = 491
$ws.Range("Q8").Value = 6982
$ws.Range("R8").Value = -4407
$ws.Range("S8").Value = -2054
$ws.Range("T8").Value = 2917
$ws.Range("U8").Value = 4669
$ws.Range("W8").Value = 7.18
$ws.Range("X8").Value = 4.62
$ws.Range("Y8").Value = 6.24
$ws.Range("Z8").Value = 2.35
$ws.Range("AA8").Value = 127.64
$ws.Range("AC8").Value = 27395
$ws.Range("AD8").Value = 9.75
$ws.Range("AE8").Value = 452362
$ws.Range("AF8").Value = 0.59
$ws.Range("AG8").Value = 2086
$ws.Range("AH8").Value = 0.78
$ws.Range("AI8").Value = 7.61
$ws.Range("D9").Value = 75046
$ws.Range("E9").Value = 5639
$ws.Range("G9").Value = 4989
$ws.Range("H9").Value = 3783
$ws.Range("I9").Value = 3188
$ws.Range("K9").Value = 140917
$ws.Range("L9").Value = 76885
$ws.Range("M9").Value = 64031
$ws.Range("N9").Value = 47410
$ws.Range("P9").Value = 491
$ws.Range("Q9").Value = 6784
$ws.Range("R9").Value = -4362
$ws.Range("S9").Value = -1896
$ws.Range("T9").Value = 2888
$ws.Range("U9").Value = 4591
$ws.Range("W9").Value = 7.51
$ws.Range("X9").Value = 5.04
$ws.Range("Y9").Value = 6.94
$ws.Range("Z9").Value = 2.71
$ws.Range("AA9").Value = 120.08
$ws.Range("AC9").Value = 32386
$ws.Range("AD9").Value = 8.24
$ws.Range("AE9").Value = 481988
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 2167
$ws.Range("AH9").Value = 0.8100000000000001
$ws.Range("AI9").Value = 6.69
